# Update column F (dSF) values on the active worksheet to reflect the
# repulled/recalculated data, per the commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    3  = -7
    4  = -6
    6  = -2
    9  = -1
    11 = -5
    12 = -1
    16 = -4
    17 = -2
    18 = 1
    20 = -2
    28 = 0
    43 = 1
    46 = -1
    47 = 0
    48 = 1
    58 = 0
    63 = -6
    64 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
